$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2785.9
$ws.Range("I70").Value = 3329.1667
$ws.Range("J70").Value = 1971
$ws.Range("K70").Value = 9987.500100000001
$ws.Range("L70").Value = 5913
$ws.Range("M70").Value = -9717.500100000001
$ws.Range("N70").Value = -6453

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2785.9
$ws.Range("I73").Value = 3329.1667
$ws.Range("J73").Value = 1971
$ws.Range("K73").Value = 9987.500100000001
$ws.Range("L73").Value = 5913
$ws.Range("M73").Value = -9051.500100000001
$ws.Range("N73").Value = -7785

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7301.4
$ws.Range("I2").Value = 4599.6
$ws.Range("K2").Value = 4599.6
$ws.Range("M2").Value = -4486.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6008.8887
$ws.Range("I32").Value = 4089.923
$ws.Range("K32").Value = 4089.923
$ws.Range("M32").Value = -3802.923

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3914.5715
$ws.Range("I45").Value = 2365
$ws.Range("J45").Value = 5076.75
$ws.Range("K45").Value = 2365
$ws.Range("L45").Value = 5076.75
$ws.Range("M45").Value = -1988
$ws.Range("N45").Value = -5830.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 743.2917
$ws.Range("I97").Value = 515.4091
$ws.Range("J97").Value = 3250
$ws.Range("K97").Value = 515.4091
$ws.Range("L97").Value = 3250
$ws.Range("M97").Value = -19.40909999999997
$ws.Range("N97").Value = -4242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 7301.4
$ws.Range("I116").Value = 4599.6
$ws.Range("K116").Value = 4599.6
$ws.Range("M116").Value = -2305.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7301.4
$ws.Range("I3").Value = 4599.6
$ws.Range("K3").Value = 4599.6
$ws.Range("M3").Value = -4485.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2285.577
$ws.Range("I94").Value = 1766.7368
$ws.Range("K94").Value = 1766.7368
$ws.Range("M94").Value = -1315.7368

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 983.6
$ws.Range("I107").Value = 640.6667
$ws.Range("K107").Value = 640.6667
$ws.Range("M107").Value = 1279.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1091.409
$ws.Range("I134").Value = 1000.55
$ws.Range("K134").Value = 3001.65
$ws.Range("M134").Value = -466.6499999999996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 999
$ws.Range("I16").Value = 999
$ws.Range("K16").Value = 999
$ws.Range("M16").Value = -712

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1867.5555
$ws.Range("I31").Value = 1381.8
$ws.Range("K31").Value = 1381.8
$ws.Range("M31").Value = -1086.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1867.5555
$ws.Range("I34").Value = 1381.8
$ws.Range("K34").Value = 1381.8
$ws.Range("M34").Value = -1179.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2184
$ws.Range("I105").Value = 2078.6667
$ws.Range("K105").Value = 2078.6667
$ws.Range("M105").Value = -331.6667000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 299.33334
$ws.Range("I107").Value = 303.2
$ws.Range("K107").Value = 303.2
$ws.Range("M107").Value = 1616.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 999
$ws.Range("I113").Value = 999
$ws.Range("K113").Value = 999
$ws.Range("M113").Value = 1171

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2573.0833
$ws.Range("I132").Value = 2827.9048
$ws.Range("K132").Value = 8483.714399999999
$ws.Range("M132").Value = -5953.714399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 46.346153
$ws.Range("I2").Value = 29.210526
$ws.Range("J2").Value = 92.85714
$ws.Range("K2").Value = 175.263156
$ws.Range("L2").Value = 557.14284
$ws.Range("M2").Value = -62.26315600000001
$ws.Range("N2").Value = -783.14284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3168.5625
$ws.Range("I69").Value = 2245.4546
$ws.Range("K69").Value = 6736.3638
$ws.Range("M69").Value = -5925.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 3168.5625
$ws.Range("I72").Value = 2245.4546
$ws.Range("K72").Value = 20209.0914
$ws.Range("M72").Value = -16153.0914

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 562.6667
$ws.Range("I92").Value = 296
$ws.Range("K92").Value = 888
$ws.Range("M92").Value = 360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2491.889
$ws.Range("J103").Value = 2867.5
$ws.Range("L103").Value = 8602.5
$ws.Range("N103").Value = -10360.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 11950
$ws.Range("J92").Value = 11950
$ws.Range("L92").Value = 11950
$ws.Range("N92").Value = -15694

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2423.2222
$ws.Range("I122").Value = 2601.125
$ws.Range("K122").Value = 7803.375
$ws.Range("M122").Value = -5353.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2949.75
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2949.75
$ws.Range("I27").Value = 300
$ws.Range("K27").Value = 300
$ws.Range("M27").Value = -193

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1684.9565
$ws.Range("I46").Value = 1068.1818
$ws.Range("K46").Value = 1068.1818
$ws.Range("M46").Value = -880.1818000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2800
$ws.Range("I68").Value = 2800
$ws.Range("K68").Value = 2800
$ws.Range("M68").Value = -2051

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2800
$ws.Range("I71").Value = 2800
$ws.Range("K71").Value = 14000
$ws.Range("M71").Value = -10256

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 34183332
$ws.Range("I3").Value = 100000000
$ws.Range("J3").Value = 1275000
$ws.Range("K3").Value = 100000000
$ws.Range("L3").Value = 1275000
$ws.Range("M3").Value = -99999886
$ws.Range("N3").Value = -1275228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 6337983
$ws.Range("I100").Value = 13940304
$ws.Range("J100").Value = 2716.1667
$ws.Range("K100").Value = 27880608
$ws.Range("L100").Value = 5432.3334
$ws.Range("M100").Value = -27880067
$ws.Range("N100").Value = -6514.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 899.5
$ws.Range("J107").Value = 799
$ws.Range("L107").Value = 2397
$ws.Range("N107").Value = -6237

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 536.8
$ws.Range("I113").Value = 296
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 888
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = 1282
$ws.Range("N113").Value = -8840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2379.6
$ws.Range("I132").Value = 2199.2778
$ws.Range("J132").Value = 4002.5
$ws.Range("K132").Value = 6597.8334
$ws.Range("L132").Value = 12007.5
$ws.Range("M132").Value = -4067.8334
$ws.Range("N132").Value = -17067.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
